$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number formatting (styles) from the last existing data row (317)
# down into the new rows (318:324) for columns A (date) and B (depth),
# reusing the same cellXf indices rather than creating new styles.
[void]$ws.Range("A317:B317").Copy()
[void]$ws.Range("A318:B324").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(318, 1).Value = 45361.55
$ws.Cells.Item(318, 2).Value = 1.83

$ws.Cells.Item(319, 1).Value = 45362.652777777781
$ws.Cells.Item(319, 2).Value = 1.82
$ws.Cells.Item(319, 4).Value = "Flynn lake ice out"

$ws.Cells.Item(320, 1).Value = 45363.5
$ws.Cells.Item(320, 2).Value = 1.82

$ws.Cells.Item(321, 1).Value = 45364.668055555558
$ws.Cells.Item(321, 2).Value = 1.82

$ws.Cells.Item(322, 1).Value = 45365.63958333333
$ws.Cells.Item(322, 2).Value = 1.8
$ws.Cells.Item(322, 4).Value = "Eagle lake ice out"

$ws.Cells.Item(323, 1).Value = 45366.557638888888
$ws.Cells.Item(323, 2).Value = 1.84

$ws.Cells.Item(324, 1).Value = 45367.664583333331
$ws.Cells.Item(324, 2).Value = 1.85

# Update the frozen-pane view and active selection to reflect the new
# bottom of the data (matches what Excel records when scrolled to the end).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 300
[void]$ws.Range("A2").Select()
[void]$ws.Range("C324").Select()
